$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Merge A1:A3 and set value/alignment to match header cells (Benchmark column spans header rows)
$ws.Range("A3").ClearContents()
$ws.Range("A1").Value = "Benchmark"
$ws.Range("A1:A3").Merge()
$ws.Range("A1:A3").HorizontalAlignment = -4108  # xlCenter
$ws.Range("A1:A3").VerticalAlignment = -4108    # xlCenter

# Adjust column widths
$ws.Columns.Item(2).ColumnWidth = 12
$ws.Columns.Item(3).ColumnWidth = 13.7109375
$ws.Columns.Item(4).ColumnWidth = 12
$ws.Columns.Item(5).ColumnWidth = 13.7109375
$ws.Columns.Item(6).ColumnWidth = 12
$ws.Columns.Item(7).ColumnWidth = 13.7109375
$ws.Columns.Item(8).ColumnWidth = 12
$ws.Columns.Item(9).ColumnWidth = 13.7109375

# Move selection
$ws.Range("D21").Select()
